$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('E21').Value = '  +3.15%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('E23').Value = '  -3.43%  '
$ws.Range('E24').Value = '  +3.26%  '
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  +5.67%  '
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('E40').Value = '  -3.51%  '
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('E42').Value = '  +3.21%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('E51').Value = '  -1.09%  '

# Price column (D) values are forced to text to avoid Excel auto-numeric conversion
# (matches original inline-string cell typing and preserves exact formatting/trailing zeros)
$ws.Range('D2').Value = "'68.189.28"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'2.642.39"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Value = "'597.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'156.19"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D12').Value = "'0.350"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'27.96"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Value = "'3.122.10"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'68.115.01"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'2.637.59"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Value = "'363.41"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'7.33"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D24').Value = "'75.38"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Value = "'9.74"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'1.06"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'2.775.57"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Value = "'555.74"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = "'8.05"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D37').Value = "'161.55"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'19.64"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D42').Value = "'0.0₆0333"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'17.80"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D46').Value = "'158.55"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D49').Value = "'0.0782"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'1.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'0.614"
$ws.Range('D51').Style = 'Normal'
